$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 177, shifting existing rows 177-262 down to 178-263.
$ws.Rows.Item(177).Insert()

# Populate the newly inserted row 177 with the new record.
$ws.Range("A177").Value = 6
$ws.Range("B177").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C177").Value = "Metropolitana"
$ws.Range("D177").Value = 44839
$ws.Range("E177").Value = 13
$ws.Range("F177").Value = 100112022
$ws.Range("G177").Value = "Arveja Verde"
$ws.Range("H177").Value = "Perfection"
$ws.Range("I177").Value = "Primera"
$ws.Range("J177").Value = 400
$ws.Range("K177").Value = 25000
$ws.Range("L177").Value = 27000
$ws.Range("M177").Value = 25850
$ws.Range("N177").Value = '$/malla 25 kilos'
$ws.Range("O177").Value = "Provincia de Huasco"
$ws.Range("P177").Value = 1034
$ws.Range("Q177").Value = 25
$ws.Range("R177").Value = "Hortaliza"
